$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.229.68"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "1.817.84"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'313.17"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.4662"
$ws.Range("E7").Value = "  +4.42%  "
$ws.Range("D8").Value = "'0.3764"
$ws.Range("E8").Value = "  +2.30%  "
$ws.Range("D9").Value = "'0.07392"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").Value = "'0.8703"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("D11").Value = "'20.59"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "1.824.39"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").Value = "'6.673"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").Value = "'5.402"
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").Value = "'92.16"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "'0.07079"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "'0.000008755"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "'14.91"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "27.262.63"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").Value = "'5.307"
$ws.Range("E22").Value = "  +2.95%  "
$ws.Range("D23").Value = "'10.93"
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").Value = "2.050.03"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "'1.939"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("D26").Value = "'151.56"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "'2.234"
$ws.Range("E27").Value = "  +2.92%  "
$ws.Range("D28").Value = "'18.55"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").Value = "'5.295"
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("D30").Value = "'117.07"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'0.08936"
$ws.Range("E31").Value = "  +1.81%  "
$ws.Range("D32").Value = "'0.7810"
$ws.Range("E32").Value = "  +5.38%  "
$ws.Range("D33").Value = "'1.178"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").Value = "'4.519"
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("D35").Value = "'2.932"
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "'1.098"
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").Value = "'0.05242"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").Value = "'7.268"
$ws.Range("E40").Value = "  +4.32%  "
$ws.Range("D41").Value = "'0.5314"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'2.367"
$ws.Range("E42").Value = "  +20.65%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "'2.882"
$ws.Range("E43").Value = "  +1.71%  "
$ws.Range("D44").Value = "'0.1688"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").Value = "'8.573"
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("D46").Value = "'0.5050"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").Value = "'10.45"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Value = "'105.44"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").Value = "'1.000"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "'1.665"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").Value = "'0.06337"
$ws.Range("E51").Value = "  +0.73%  "
